$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 ("테스트" player). This shifts the
# existing rows 4-7 (페이커, 차두리, 손흥민, 박지성) down to rows 5-8,
# carrying their styles (s="1") with them.
$ws.Rows.Item(4).Insert()

# Fill in the newly inserted row 4: id=3, name="테스트", all stats = 0.
$ws.Cells.Item(4,1).Value = 3
$ws.Cells.Item(4,2).Value = "테스트"
$ws.Cells.Item(4,3).Value = 0
$ws.Cells.Item(4,4).Value = 0
$ws.Cells.Item(4,5).Value = 0
$ws.Cells.Item(4,6).Value = 0
$ws.Cells.Item(4,7).Value = 0

# Row 6 (차두리) stat values changed from 80/80/80/50/50 to 50/50/50/50/50.
$ws.Cells.Item(6,3).Value = 50
$ws.Cells.Item(6,4).Value = 50
$ws.Cells.Item(6,5).Value = 50
$ws.Cells.Item(6,6).Value = 50
$ws.Cells.Item(6,7).Value = 50

# Append a brand new row 9 ("테스트2" player): id=8, all stats = 50.
# Use the style from the row above (row 8, 박지성) so the new row
# matches the existing formatting (s="1").
$ws.Range("A8:G8").Copy()
$ws.Range("A9:G9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(9,1).Value = 8
$ws.Cells.Item(9,2).Value = "테스트2"
$ws.Cells.Item(9,3).Value = 50
$ws.Cells.Item(9,4).Value = 50
$ws.Cells.Item(9,5).Value = 50
$ws.Cells.Item(9,6).Value = 50
$ws.Cells.Item(9,7).Value = 50
